$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item(1)

$src.Range("A1:J1").Select()

$existing = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $existing)
$ws.Name = "new_customers"

$headerRange = $src.Range("A1:I1")
$dataRange = $src.Range("A961:I999")

$headerRange.Copy()
$ws.Range("A1:I1").PasteSpecial()

$dataRange.Copy()
$ws.Range("A2:I40").PasteSpecial()
